# Generate Report for Handoff
# Replace the old handoff id (ddf358cc-870d-4c68-994d-1adea9ae8ca0) with the
# new one (2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed) and bump the handoff
# timestamps across the three report sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.md"
$ws.Range("D2").Value = "2016-03-21 03:40:09"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add(
    $ws.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/4f77eb76e20774a5b6582e1a12c7ede7286b9c4c/e2e/ddf358cc-870d-4c68-994d-1adea9ae8ca0.md",
    "",
    "",
    "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.md"
)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.md"
$ws.Range("D2").Value = "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.e16630ceba8d47df1b7ae033d6fa9364a20155a6.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-21 03:40:00"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()

$ws.Hyperlinks.Add(
    $ws.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/4f77eb76e20774a5b6582e1a12c7ede7286b9c4c/e2e/ddf358cc-870d-4c68-994d-1adea9ae8ca0.md",
    "",
    "",
    "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.md"
)
$ws.Hyperlinks.Add(
    $ws.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/292b8f4430ff033702462c864dedd6448334b35f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/ddf358cc-870d-4c68-994d-1adea9ae8ca0.24b4cf000011bc438b1674c98c5e5caff0586c0b.zh-cn.xlf",
    "",
    "",
    "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.e16630ceba8d47df1b7ae033d6fa9364a20155a6.zh-cn.xlf"
)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.md"
$ws.Range("D2").Value = "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.e16630ceba8d47df1b7ae033d6fa9364a20155a6.de-de.xlf"
$ws.Range("E2").Value = "2016-03-21 03:40:09"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()

$ws.Hyperlinks.Add(
    $ws.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/4f77eb76e20774a5b6582e1a12c7ede7286b9c4c/e2e/ddf358cc-870d-4c68-994d-1adea9ae8ca0.md",
    "",
    "",
    "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.md"
)
$ws.Hyperlinks.Add(
    $ws.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/394b2c3b7b25e626829e1dd7c6913f7fc9912687/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/ddf358cc-870d-4c68-994d-1adea9ae8ca0.24b4cf000011bc438b1674c98c5e5caff0586c0b.de-de.xlf",
    "",
    "",
    "2fe20f3b-23d5-4b89-9f4a-5f0d73a70bed.e16630ceba8d47df1b7ae033d6fa9364a20155a6.de-de.xlf"
)
